$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 - this shifts the existing rows 12..33 down to 13..34
$ws.Rows.Item(12).Insert()

# Fill in the new weekly record in row 12
$ws.Cells.Item(12, 1).Value2 = 11
$ws.Cells.Item(12, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value2 = "Bíobío"
$ws.Cells.Item(12, 4).Value2 = 44797
$ws.Cells.Item(12, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value2 = 8
$ws.Cells.Item(12, 6).Value2 = 100112022
$ws.Cells.Item(12, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(12, 8).Value2 = "Perfection"
$ws.Cells.Item(12, 9).Value2 = "Primera"
$ws.Cells.Item(12, 10).Value2 = 100
$ws.Cells.Item(12, 11).Value2 = 37000
$ws.Cells.Item(12, 12).Value2 = 38000
$ws.Cells.Item(12, 13).Value2 = 37500
$ws.Cells.Item(12, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value2 = "Provincia de Huasco"
$ws.Cells.Item(12, 16).Value2 = 1500
$ws.Cells.Item(12, 17).Value2 = 25
$ws.Cells.Item(12, 18).Value2 = "Hortaliza"
